$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '29.072.80'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '1.822.55'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '''233.77'
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').Value = '''0.5986'
$ws.Range('E6').Value = '  -4.09%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''0.06938'
$ws.Range('E8').Value = '  -5.90%  '
$ws.Range('D9').Value = '''0.2745'
$ws.Range('E9').Value = '  -4.78%  '
$ws.Range('D10').Value = '''23.19'
$ws.Range('E10').Value = '  -6.15%  '
$ws.Range('D11').Value = '''0.07586'
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('D12').Value = '1.825.92'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '''4.729'
$ws.Range('E13').Value = '  -4.37%  '
$ws.Range('D14').Value = '''0.6233'
$ws.Range('E14').Value = '  -5.70%  '
$ws.Range('D15').Value = '''0.000009734'
$ws.Range('E15').Value = '  -7.67%  '
$ws.Range('D16').Value = '''77.20'
$ws.Range('E16').Value = '  -4.91%  '
$ws.Range('D17').Value = '28.708.21'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').Value = '''5.564'
$ws.Range('E18').Value = '  -10.83%  '
$ws.Range('D19').Value = '''215.32'
$ws.Range('E19').Value = '  -7.52%  '
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = '''11.50'
$ws.Range('E21').Value = '  -5.91%  '
$ws.Range('D22').Value = '''6.843'
$ws.Range('E22').Value = '  -6.18%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '''156.15'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').Value = '''7.906'
$ws.Range('E25').Value = '  -6.04%  '
$ws.Range('D26').Value = '''0.1281'
$ws.Range('E26').Value = '  -3.95%  '
$ws.Range('D27').Value = '''16.40'
$ws.Range('E27').Value = '  -4.84%  '
$ws.Range('D28').Value = '''0.06451'
$ws.Range('E28').Value = '  -9.06%  '
$ws.Range('D29').Value = '''1.426'
$ws.Range('E29').Value = '  -3.99%  '
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('D31').Value = '''3.825'
$ws.Range('E31').Value = '  -4.80%  '
$ws.Range('D32').Value = '''3.742'
$ws.Range('E32').Value = '  -7.04%  '
$ws.Range('D33').Value = '''1.719'
$ws.Range('E33').Value = '  -4.38%  '
$ws.Range('D34').Value = '''1.086'
$ws.Range('E34').Value = '  -5.40%  '
$ws.Range('D35').Value = '''0.6444'
$ws.Range('E35').Value = '  -7.17%  '
$ws.Range('D36').Value = '''2.533'
$ws.Range('E36').Value = '  -2.02%  '
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').Value = '''0.01743'
$ws.Range('E38').Value = '  -4.42%  '
$ws.Range('D39').Value = '''6.513'
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('D40').Value = '1.144.54'
$ws.Range('E40').Value = '  -7.20%  '
$ws.Range('D41').Value = '''0.8804'
$ws.Range('E41').Value = '  -7.50%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '1.973.34'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').Value = '''100.15'
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('D45').Value = '''61.38'
$ws.Range('E45').Value = '  -5.77%  '
$ws.Range('D46').Value = '''0.00000000112'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('D47').Value = '''1.607'
$ws.Range('E47').Value = '  -4.14%  '
$ws.Range('D48').Value = '''8.427'
$ws.Range('E48').Value = '  -5.26%  '
$ws.Range('D49').Value = '''0.05500'
$ws.Range('E49').Value = '  -2.59%  '
$ws.Range('D50').Value = '''0.4532'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').Value = '''6.399'
$ws.Range('E51').Value = '  -7.61%  '
